# Apply changes described by the diff:
# 1. Rename sheet "ICU" to "Ventilation"
# 2. Update selections (active cell) on each sheet
# 3. Make "Ventilation" the active sheet (tabSelected / activeTab)

$wb = $excel.ActiveWorkbook

$wsSeverity = $wb.Worksheets.Item("Severity")
$wsHospit = $wb.Worksheets.Item("Hospit")
$wsICU = $wb.Worksheets.Item("ICU")

# Rename ICU sheet to Ventilation
$wsICU.Name = "Ventilation"

# Update selection on Severity sheet
$wsSeverity.Activate()
$wsSeverity.Range("D13").Select()

# Update selection on Hospit sheet
$wsHospit.Activate()
$wsHospit.Range("B39").Select()

# Update selection on Ventilation sheet (formerly ICU) and make it active
$wsVentilation = $wb.Worksheets.Item("Ventilation")
$wsVentilation.Activate()
$wsVentilation.Range("H38").Select()
